$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "58.071.06"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "2.349.75"
$ws.Range("E3").Value = "  +1.30%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "545.56"
$ws.Range("E5").Value = "  +5.81%  "
$ws.Range("D6").Value = "134.93"
$ws.Range("E6").Value = "  +2.36%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "0.538"
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "2.348.01"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("E11").Value = "  +1.12%  "
$ws.Range("D12").Value = "5.41"
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("E13").Value = "  +6.44%  "
$ws.Range("D14").Value = "2.767.09"
$ws.Range("E14").Value = "  +1.24%  "
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "58.045.25"
$ws.Range("E16").Value = "  +2.53%  "
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "2.352.43"
$ws.Range("E18").Value = "  +1.12%  "
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("D20").Value = "334.63"
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("E21").Value = "  +1.68%  "
$ws.Range("D22").Value = "6.70"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").Value = "61.78"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "8.44"
$ws.Range("E27").Value = "  -2.01%  "
$ws.Range("D28").Value = "1.42"
$ws.Range("E28").Value = "  +8.13%  "
$ws.Range("E29").Value = "  +5.41%  "
$ws.Range("D30").Value = "170.29"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").Value = "0.0₃0729"
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("E33").Value = "  +17.77%  "
$ws.Range("D34").Value = "18.43"
$ws.Range("E34").Value = "  +0.78%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  +6.02%  "
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("E39").Value = "  +4.55%  "
$ws.Range("D41").Value = "147.90"
$ws.Range("D42").Value = "0.378"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").Value = "288.02"
$ws.Range("E43").Value = "  +4.39%  "
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").Value = "19.26"
$ws.Range("E45").Value = "  +5.35%  "
$ws.Range("E47").Value = "  +1.93%  "
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("E49").Value = "  +1.51%  "
$ws.Range("D50").Value = "17.57"
$ws.Range("E50").Value = "  +2.96%  "
$ws.Range("E51").Value = "  +1.20%  "
